$wb = $excel.ActiveWorkbook

# ===== Sheet 1: ALC =====
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 798
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H3").Value = 20657
$ws.Range("J3").Value = 20657
$ws.Range("L3").Value = 20657
$ws.Range("N3").Value = -20885
$ws.Range("H8").Value = 254
$ws.Range("I8").Value = 254
$ws.Range("K8").Value = 762
$ws.Range("M8").Value = -623
$ws.Range("H9").Value = 260
$ws.Range("I9").Value = 274.5
$ws.Range("J9").Value = 202
$ws.Range("K9").Value = 274.5
$ws.Range("L9").Value = 202
$ws.Range("M9").Value = -105.5
$ws.Range("N9").Value = -540
$ws.Range("H11").Value = 362.83334
$ws.Range("I11").Value = 362.83334
$ws.Range("K11").Value = 362.83334
$ws.Range("M11").Value = -222.83334
$ws.Range("H18").Value = 807.5
$ws.Range("I18").Value = 807.5
$ws.Range("K18").Value = 807.5
$ws.Range("M18").Value = -523.5
$ws.Range("H28").Value = 3429.1667
$ws.Range("I28").Value = 2514.2
$ws.Range("J28").Value = 8004
$ws.Range("K28").Value = 2514.2
$ws.Range("L28").Value = 8004
$ws.Range("M28").Value = -2029.2
$ws.Range("N28").Value = -8974
$ws.Range("H33").Value = 761.3913
$ws.Range("I33").Value = 212.26666
$ws.Range("K33").Value = 212.26666
$ws.Range("M33").Value = 16.73334
$ws.Range("H40").Value = 5724.909
$ws.Range("I40").Value = 5634.5
$ws.Range("K40").Value = 5634.5
$ws.Range("M40").Value = -5459.5
$ws.Range("H43").Value = 1995.5
$ws.Range("I43").Value = 999
$ws.Range("J43").Value = 2992
$ws.Range("K43").Value = 999
$ws.Range("L43").Value = 2992
$ws.Range("M43").Value = -930
$ws.Range("N43").Value = -3130
$ws.Range("H48").Value = 4798.5
$ws.Range("J48").Value = 4798.5
$ws.Range("L48").Value = 14395.5
$ws.Range("N48").Value = -14979.5
$ws.Range("H56").Value = 4798.5
$ws.Range("J56").Value = 4798.5
$ws.Range("L56").Value = 14395.5
$ws.Range("N56").Value = -15463.5
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880
$ws.Range("H70").Value = 7747.25
$ws.Range("J70").Value = 10488.546
$ws.Range("L70").Value = 31465.638
$ws.Range("N70").Value = -32005.638
$ws.Range("H73").Value = 7747.25
$ws.Range("J73").Value = 10488.546
$ws.Range("L73").Value = 31465.638
$ws.Range("N73").Value = -33337.638
$ws.Range("H80").Value = 916.2
$ws.Range("I80").Value = 791
$ws.Range("K80").Value = 2373
$ws.Range("M80").Value = -1375
$ws.Range("H83").Value = 916.2
$ws.Range("I83").Value = 791
$ws.Range("K83").Value = 7119
$ws.Range("M83").Value = -2127
$ws.Range("H88").Value = 6340.7144
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 6980.8335
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 6980.8335
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -7792.8335
$ws.Range("H91").Value = 6340.7144
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 6980.8335
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 6980.8335
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -9788.8335
$ws.Range("H97").Value = 958.4286
$ws.Range("I97").Value = 999
$ws.Range("K97").Value = 2997
$ws.Range("M97").Value = -2501
$ws.Range("H100").Value = 6482.1665
$ws.Range("J100").Value = 5900
$ws.Range("L100").Value = 5900
$ws.Range("N100").Value = -6982
$ws.Range("H102").Value = 20657
$ws.Range("J102").Value = 20657
$ws.Range("L102").Value = 20657
$ws.Range("N102").Value = -27147
$ws.Range("H107").Value = 1391.2858
$ws.Range("I107").Value = 1139
$ws.Range("K107").Value = 1139
$ws.Range("M107").Value = 781
$ws.Range("H116").Value = 9500
$ws.Range("I116").Value = 9500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 9500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -6058
$ws.Range("N116").ClearContents()
$ws.Range("H127").Value = 960
$ws.Range("I127").Value = 787.44446
$ws.Range("K127").Value = 2362.33338
$ws.Range("M127").Value = 2597.66662
$ws.Range("H129").Value = 1673.875
$ws.Range("I129").Value = 1673.875
$ws.Range("K129").Value = 5021.625
$ws.Range("M129").Value = -21.625
$ws.Range("H132").Value = 2246.982
$ws.Range("I132").Value = 2105.2454
$ws.Range("K132").Value = 6315.736199999999
$ws.Range("M132").Value = -3785.736199999999
$ws.Range("H135").Value = 5532.96
$ws.Range("I135").Value = 1196.591
$ws.Range("K135").Value = 10769.319
$ws.Range("M135").Value = -8234.319
$ws.Range("H137").Value = 3102.1667
$ws.Range("I137").Value = 3377.5
$ws.Range("J137").Value = 2551.5
$ws.Range("K137").Value = 10132.5
$ws.Range("L137").Value = 7654.5
$ws.Range("M137").Value = -7582.5
$ws.Range("N137").Value = -12754.5
$ws.Range("H138").Value = 2809.0981
$ws.Range("J138").Value = 3097.762
$ws.Range("L138").Value = 9293.286
$ws.Range("N138").Value = -19573.286
$ws.Range("H140").Value = 114575.8
$ws.Range("J140").Value = 114575.8
$ws.Range("L140").Value = 114575.8
$ws.Range("N140").Value = -124935.8

# ===== Sheet 2: ARM =====
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 3178.1667
$ws.Range("J2").Value = 2533
$ws.Range("L2").Value = 2533
$ws.Range("N2").Value = -2759
$ws.Range("H28").Value = 21333
$ws.Range("I28").Value = 21333
$ws.Range("K28").Value = 21333
$ws.Range("M28").Value = -21141
$ws.Range("H32").Value = 7368.877
$ws.Range("I32").Value = 889.8
$ws.Range("K32").Value = 889.8
$ws.Range("M32").Value = -602.8
$ws.Range("H45").Value = 10271.593
$ws.Range("I45").Value = 12237.429
$ws.Range("K45").Value = 12237.429
$ws.Range("M45").Value = -11860.429
$ws.Range("H61").Value = 2621.111
$ws.Range("I61").Value = 2485.2122
$ws.Range("K61").Value = 2485.2122
$ws.Range("M61").Value = -2273.2122
$ws.Range("H74").Value = 2533.353
$ws.Range("I74").Value = 2559.5625
$ws.Range("J74").Value = 2114
$ws.Range("K74").Value = 2559.5625
$ws.Range("L74").Value = 2114
$ws.Range("M74").Value = -1685.5625
$ws.Range("N74").Value = -3862
$ws.Range("H77").Value = 2533.353
$ws.Range("I77").Value = 2559.5625
$ws.Range("J77").Value = 2114
$ws.Range("K77").Value = 12797.8125
$ws.Range("L77").Value = 10570
$ws.Range("M77").Value = -8429.8125
$ws.Range("N77").Value = -19306
$ws.Range("H97").Value = 1004.15
$ws.Range("I97").Value = 616.64703
$ws.Range("K97").Value = 616.64703
$ws.Range("M97").Value = -120.64703
$ws.Range("H99").Value = 21333
$ws.Range("I99").Value = 21333
$ws.Range("K99").Value = 21333
$ws.Range("M99").Value = -18338
$ws.Range("H110").Value = 2852.75
$ws.Range("J110").Value = 3085
$ws.Range("L110").Value = 3085
$ws.Range("N110").Value = -7175
$ws.Range("H116").Value = 3178.1667
$ws.Range("J116").Value = 2533
$ws.Range("L116").Value = 2533
$ws.Range("N116").Value = -7121
$ws.Range("H132").Value = 2357
$ws.Range("I132").Value = 2749.5
$ws.Range("K132").Value = 8248.5
$ws.Range("M132").Value = -5718.5
$ws.Range("H136").Value = 2621.111
$ws.Range("I136").Value = 2485.2122
$ws.Range("K136").Value = 7455.6366
$ws.Range("M136").Value = -4905.6366

# ===== Sheet 3: BSM =====
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 3178.1667
$ws.Range("J3").Value = 2533
$ws.Range("L3").Value = 2533
$ws.Range("N3").Value = -2761
$ws.Range("H57").Value = 80569.664
$ws.Range("J57").Value = 80000
$ws.Range("L57").Value = 80000
$ws.Range("N57").Value = -81440
$ws.Range("H80").Value = 595.1429
$ws.Range("I80").Value = 133.33333
$ws.Range("K80").Value = 133.33333
$ws.Range("M80").Value = 864.6666700000001
$ws.Range("H83").Value = 595.1429
$ws.Range("I83").Value = 133.33333
$ws.Range("K83").Value = 666.6666499999999
$ws.Range("M83").Value = 4325.33335
$ws.Range("H86").Value = 3190.38
$ws.Range("I86").Value = 2765.15
$ws.Range("K86").Value = 2765.15
$ws.Range("M86").Value = -1642.15
$ws.Range("H89").Value = 3190.38
$ws.Range("I89").Value = 2765.15
$ws.Range("K89").Value = 13825.75
$ws.Range("M89").Value = -8209.75
$ws.Range("H94").Value = 1704.9474
$ws.Range("I94").Value = 1219.2
$ws.Range("J94").Value = 2021.7391
$ws.Range("K94").Value = 1219.2
$ws.Range("L94").Value = 2021.7391
$ws.Range("M94").Value = -768.2
$ws.Range("N94").Value = -2923.7391
$ws.Range("H134").Value = 2522.625
$ws.Range("I134").Value = 2522.625
$ws.Range("K134").Value = 7567.875
$ws.Range("M134").Value = -5032.875
$ws.Range("H136").Value = 80569.664
$ws.Range("J136").Value = 80000
$ws.Range("L136").Value = 80000
$ws.Range("N136").Value = -90200
$ws.Range("H137").Value = 77774.5
$ws.Range("J137").Value = 77774.5
$ws.Range("L137").Value = 77774.5
$ws.Range("N137").Value = -87974.5

# ===== Sheet 4: CRP =====
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 5707.6665
$ws.Range("I31").Value = 9221.111
$ws.Range("J31").Value = 3072.5833
$ws.Range("K31").Value = 9221.111
$ws.Range("L31").Value = 3072.5833
$ws.Range("M31").Value = -8926.111
$ws.Range("N31").Value = -3662.5833
$ws.Range("H34").Value = 5707.6665
$ws.Range("I34").Value = 9221.111
$ws.Range("J34").Value = 3072.5833
$ws.Range("K34").Value = 9221.111
$ws.Range("L34").Value = 3072.5833
$ws.Range("M34").Value = -9019.111
$ws.Range("N34").Value = -3476.5833
$ws.Range("H43").Value = 66006.836
$ws.Range("J43").Value = 66006.836
$ws.Range("L43").Value = 66006.836
$ws.Range("N43").Value = -66374.836
$ws.Range("H58").Value = 7862.5454
$ws.Range("I58").Value = 7862.5454
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 7862.5454
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -7659.5454
$ws.Range("N58").ClearContents()
$ws.Range("H99").Value = 8035.1333
$ws.Range("I99").Value = 6104.1
$ws.Range("J99").Value = 11897.2
$ws.Range("K99").Value = 6104.1
$ws.Range("L99").Value = 11897.2
$ws.Range("M99").Value = -4606.1
$ws.Range("N99").Value = -14893.2
$ws.Range("H101").Value = 66006.836
$ws.Range("J101").Value = 66006.836
$ws.Range("L101").Value = 66006.836
$ws.Range("N101").Value = -72496.836
$ws.Range("H124").Value = 72762.336
$ws.Range("J124").Value = 68995.5
$ws.Range("L124").Value = 68995.5
$ws.Range("N124").Value = -73905.5
$ws.Range("H126").Value = 8035.1333
$ws.Range("I126").Value = 6104.1
$ws.Range("J126").Value = 11897.2
$ws.Range("K126").Value = 18312.3
$ws.Range("L126").Value = 35691.60000000001
$ws.Range("M126").Value = -15842.3
$ws.Range("N126").Value = -40631.60000000001
$ws.Range("H132").Value = 13703
$ws.Range("I132").Value = 15194.6
$ws.Range("K132").Value = 45583.8
$ws.Range("M132").Value = -43053.8
$ws.Range("H134").Value = 2814.889
$ws.Range("I134").Value = 1487.2
$ws.Range("J134").Value = 4474.5
$ws.Range("K134").Value = 4461.6
$ws.Range("L134").Value = 13423.5
$ws.Range("M134").Value = -1926.6
$ws.Range("N134").Value = -18493.5
$ws.Range("H136").Value = 7862.5454
$ws.Range("I136").Value = 7862.5454
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 23587.6362
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -21037.6362
$ws.Range("N136").ClearContents()
$ws.Range("H140").Value = 87997
$ws.Range("J140").Value = 87997
$ws.Range("L140").Value = 87997
$ws.Range("N140").Value = -98357
$ws.Range("I141").Value = 42000
$ws.Range("J141").Value = 77856.43
$ws.Range("K141").Value = 42000
$ws.Range("L141").Value = 77856.43
$ws.Range("M141").Value = -36820
$ws.Range("N141").Value = -88216.43

# ===== Sheet 5: CUL =====
$ws = $wb.Worksheets.Item(5)
$ws.Range("H7").Value = 1277.2222
$ws.Range("I7").Value = 1515
$ws.Range("K7").Value = 4545
$ws.Range("M7").Value = -4433
$ws.Range("H11").Value = 479
$ws.Range("I11").Value = 98
$ws.Range("K11").Value = 294
$ws.Range("M11").Value = -154
$ws.Range("H12").Value = 614.6667
$ws.Range("I12").Value = 558.5
$ws.Range("J12").Value = 649.2308
$ws.Range("K12").Value = 1675.5
$ws.Range("L12").Value = 1947.6924
$ws.Range("M12").Value = -1502.5
$ws.Range("N12").Value = -2293.6924
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H34").Value = 1565.75
$ws.Range("I34").Value = 348.75
$ws.Range("J34").Value = 2174.25
$ws.Range("K34").Value = 1046.25
$ws.Range("L34").Value = 6522.75
$ws.Range("M34").Value = -962.25
$ws.Range("N34").Value = -6690.75
$ws.Range("H39").Value = 8059
$ws.Range("J39").Value = 9125
$ws.Range("L39").Value = 27375
$ws.Range("N39").Value = -27963
$ws.Range("H55").Value = 4373.1
$ws.Range("J55").Value = 4373.1
$ws.Range("L55").Value = 13119.3
$ws.Range("N55").Value = -13473.3
$ws.Range("H58").Value = 4599.6665
$ws.Range("I58").Value = 4599.6665
$ws.Range("K58").Value = 13798.9995
$ws.Range("M58").Value = -13670.9995
$ws.Range("H122").Value = 2677.7646
$ws.Range("I122").Value = 628.25
$ws.Range("J122").Value = 4499.5557
$ws.Range("K122").Value = 5654.25
$ws.Range("L122").Value = 40496.0013
$ws.Range("M122").Value = -3204.25
$ws.Range("N122").Value = -45396.0013
$ws.Range("H131").Value = 2273.3125
$ws.Range("I131").Value = 1570
$ws.Range("K131").Value = 4710
$ws.Range("M131").Value = 330

# ===== Sheet 6: GSM =====
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 5929.3335
$ws.Range("I70").Value = 5929.3335
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5929.3335
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -5659.3335
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 5929.3335
$ws.Range("I73").Value = 5929.3335
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5929.3335
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -4993.3335
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 13375
$ws.Range("I80").Value = 8750
$ws.Range("J80").Value = 18000
$ws.Range("K80").Value = 8750
$ws.Range("L80").Value = 18000
$ws.Range("M80").Value = -7752
$ws.Range("N80").Value = -19996
$ws.Range("H83").Value = 13375
$ws.Range("I83").Value = 8750
$ws.Range("J83").Value = 18000
$ws.Range("K83").Value = 43750
$ws.Range("L83").Value = 90000
$ws.Range("M83").Value = -38758
$ws.Range("N83").Value = -99984
$ws.Range("H102").Value = 2382.35
$ws.Range("I102").Value = 2417.5264
$ws.Range("K102").Value = 2417.5264
$ws.Range("M102").Value = -795.5264000000002
$ws.Range("H126").Value = 1783.3334
$ws.Range("J126").Value = 1783.3334
$ws.Range("L126").Value = 5350.0002
$ws.Range("N126").Value = -10290.0002
$ws.Range("H132").Value = 3889.5557
$ws.Range("I132").Value = 3918.8635
$ws.Range("J132").Value = 3760.6
$ws.Range("K132").Value = 11756.5905
$ws.Range("L132").Value = 11281.8
$ws.Range("M132").Value = -9226.5905
$ws.Range("N132").Value = -16341.8

# ===== Sheet 7: LTW =====
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 1938.6
$ws.Range("I40").Value = 1964.5555
$ws.Range("K40").Value = 1964.5555
$ws.Range("M40").Value = -1828.5555
$ws.Range("H46").Value = 4293.6665
$ws.Range("I46").Value = 5003.263
$ws.Range("K46").Value = 5003.263
$ws.Range("M46").Value = -4815.263
$ws.Range("H61").Value = 3097.5715
$ws.Range("I61").Value = 2988.8333
$ws.Range("K61").Value = 2988.8333
$ws.Range("M61").Value = -2786.8333
$ws.Range("H113").Value = 3097.5715
$ws.Range("I113").Value = 2988.8333
$ws.Range("K113").Value = 2988.8333
$ws.Range("M113").Value = -818.8332999999998
$ws.Range("H132").Value = 2037.439
$ws.Range("I132").Value = 1767.0968
$ws.Range("J132").Value = 2875.5
$ws.Range("K132").Value = 5301.2904
$ws.Range("L132").Value = 8626.5
$ws.Range("M132").Value = -2771.2904
$ws.Range("N132").Value = -13686.5
$ws.Range("H140").Value = 117642.25
$ws.Range("J140").Value = 117642.25
$ws.Range("L140").Value = 117642.25
$ws.Range("N140").Value = -128002.25
$ws.Range("H141").Value = 87093.8
$ws.Range("J141").Value = 87326.445
$ws.Range("L141").Value = 87326.445
$ws.Range("N141").Value = -97686.445

# ===== Sheet 8: WVR =====
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 10611.5
$ws.Range("I62").Value = 12482
$ws.Range("K62").Value = 12482
$ws.Range("M62").Value = -11858
$ws.Range("H65").Value = 10611.5
$ws.Range("I65").Value = 12482
$ws.Range("K65").Value = 62410
$ws.Range("M65").Value = -59290
$ws.Range("H81").Value = 750
$ws.Range("I81").Value = 750
$ws.Range("K81").Value = 1500
$ws.Range("M81").Value = -439
$ws.Range("H84").Value = 750
$ws.Range("I84").Value = 750
$ws.Range("K84").Value = 7500
$ws.Range("M84").Value = -2196
$ws.Range("H92").Value = 94500
$ws.Range("I92").Value = 94500
$ws.Range("K92").Value = 94500
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -92004
$ws.Range("H122").Value = 1920.871
$ws.Range("I122").Value = 1810.1923
$ws.Range("J122").Value = 2496.4
$ws.Range("K122").Value = 5430.5769
$ws.Range("L122").Value = 7489.200000000001
$ws.Range("M122").Value = -2980.5769
$ws.Range("N122").Value = -12389.2
$ws.Range("H126").Value = 1921.2307
$ws.Range("I126").Value = 1921.2307
$ws.Range("K126").Value = 5763.6921
$ws.Range("M126").Value = -3293.6921
$ws.Range("H136").Value = 2693.4
$ws.Range("I136").Value = 2693.4
$ws.Range("K136").Value = 8080.200000000001
$ws.Range("M136").Value = -5530.200000000001
$ws.Range("H138").Value = 104999.5
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 104999.5
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 104999.5
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -115279.5
